# Update the LinkedIn carousel draft: swap the NTPC/GAIL joint-venture
# article for the Kakinada green-ammonia article across all 6 slides.

$p = $ppt.ActivePresentation

$title = "Kakinada: Andhra set to launch India’s first large-scale green ammonia project on Jan 17 - NewsMeter"

# Slide 1
$s = $p.Slides.Item(1)
$s.Shapes.Item(1).TextFrame.TextRange.Runs(1).Text = $title
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$tr.Paragraphs(2).Runs(1).Text = "Andhra Pradesh is set to launch India's first large-scale green ammonia project."
$tr.Paragraphs(3).Runs(1).Text = "The project is scheduled to be inaugurated on January 17."

# Slide 2
$s = $p.Slides.Item(2)
$s.Shapes.Item(1).TextFrame.TextRange.Runs(1).Text = $title
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$tr.Paragraphs(2).Runs(1).Text = "The project will be located in Kakinada, Andhra Pradesh."
$tr.Paragraphs(3).Runs(1).Text = "It aims to contribute to sustainable energy solutions in India."

# Slide 3
$s = $p.Slides.Item(3)
$s.Shapes.Item(1).TextFrame.TextRange.Runs(1).Text = $title
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$tr.Paragraphs(2).Runs(1).Text = "This initiative marks a significant step in India's green energy transition."
$tr.Paragraphs(3).Runs(1).Text = "The project is expected to enhance the production of green ammonia in the region."

# Slide 4
$s = $p.Slides.Item(4)
$s.Shapes.Item(1).TextFrame.TextRange.Runs(1).Text = $title
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$tr.Paragraphs(2).Runs(1).Text = "The green ammonia project is part of a broader strategy to reduce carbon emissions."
$tr.Paragraphs(3).Runs(1).Text = "It aligns with global trends towards sustainable energy sources."

# Slide 5
$s = $p.Slides.Item(5)
$s.Shapes.Item(1).TextFrame.TextRange.Runs(1).Text = $title
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$tr.Paragraphs(2).Runs(1).Text = "The project is anticipated to create job opportunities in the local area."
$tr.Paragraphs(3).Runs(1).Text = "It will also support the development of related industries."

# Slide 6
$s = $p.Slides.Item(6)
$s.Shapes.Item(1).TextFrame.TextRange.Runs(1).Text = $title
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$tr.Paragraphs(2).Runs(1).Text = "The launch of this project is expected to attract investments in green technology."
$tr.Paragraphs(3).Runs(1).Text = "It represents a collaboration between government and private sectors."
